$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.479.38'
$ws.Range('E2').Value = '  +2.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.840.37'
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('E4').Value = '  +1.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.78'
$ws.Range('E5').Value = '  +1.78%  '
$ws.Range('E6').Value = '  +1.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4751'
$ws.Range('E7').Value = '  +1.92%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3697'
$ws.Range('E8').Value = '  +1.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07471'
$ws.Range('E9').Value = '  +1.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8858'
$ws.Range('E10').Value = '  +2.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.47'
$ws.Range('E11').Value = '  +0.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.873.19'
$ws.Range('E12').Value = '  +6.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07347'
$ws.Range('E13').Value = '  +3.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.453'
$ws.Range('E14').Value = '  +1.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.34'
$ws.Range('E15').Value = '  +1.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.589'
$ws.Range('E16').Value = '  +1.29%  '
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008830'
$ws.Range('E18').Value = '  +1.63%  '
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.82'
$ws.Range('E20').Value = '  +1.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.501.49'
$ws.Range('E21').Value = '  +2.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.332'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.70'
$ws.Range('E23').Value = '  +0.84%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.097.16'
$ws.Range('E24').Value = '  +4.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.903'
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.40'
$ws.Range('E26').Value = '  +1.32%  '
$ws.Range('E27').Value = '  +2.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.156'
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.253'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '118.09'
$ws.Range('E30').Value = '  +2.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09007'
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7569'
$ws.Range('E32').Value = '  +0.53%  '
$ws.Range('E33').Value = '  +2.27%  '
$ws.Range('E34').Value = '  +1.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.958'
$ws.Range('E35').Value = '  +1.55%  '
$ws.Range('E36').Value = '  +1.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.107'
$ws.Range('E37').Value = '  +2.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01960'
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.003'
$ws.Range('E40').Value = '  +0.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.336'
$ws.Range('E41').Value = '  +2.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.397'
$ws.Range('E42').Value = '  +5.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5341'
$ws.Range('E43').Value = '  +0.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1662'
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.532'
$ws.Range('E45').Value = '  +1.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4922'
$ws.Range('E46').Value = '  +1.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.60'
$ws.Range('E47').Value = '  +2.51%  '
$ws.Range('E48').Value = '  +1.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.81'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.681'
$ws.Range('E50').Value = '  +1.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06322'
$ws.Range('E51').Value = '  +0.50%  '
